$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 66; $r++) {
    $ws.Range("C$r").Value = 45208
}
